$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '66.240.23'
Set-TextValue "E2" '  -0.76%  '
Set-TextValue "D3" '3.493.67'
Set-TextValue "E3" '  -0.16%  '
Set-TextValue "E4" '  -0.05%  '
Set-TextValue "D5" '604.65'
Set-TextValue "E5" '  +0.65%  '
Set-TextValue "D6" '144.32'
Set-TextValue "E6" '  -2.21%  '
Set-TextValue "D7" '3.491.56'
Set-TextValue "E7" '  -0.26%  '
Set-TextValue "E8" '  -0.03%  '
Set-TextValue "E9" '  -0.57%  '
Set-TextValue "D10" '8.06'
Set-TextValue "E10" '  +2.37%  '
Set-TextValue "D11" '0.136'
Set-TextValue "E11" '  -4.39%  '
Set-TextValue "E12" '  -2.32%  '
Set-TextValue "D13" '4.086.78'
Set-TextValue "E13" '  -0.07%  '
Set-TextValue "D14" '0.0000204'
Set-TextValue "E14" '  -4.31%  '
Set-TextValue "D15" '30.40'
Set-TextValue "E15" '  -2.54%  '
Set-TextValue "D16" '3.492.14'
Set-TextValue "E16" '  -0.43%  '
Set-TextValue "D17" '66.287.06'
Set-TextValue "E17" '  -0.71%  '
Set-TextValue "D19" '10.70'
Set-TextValue "E19" '  +2.47%  '
Set-TextValue "D20" '6.16'
Set-TextValue "E20" '  -3.52%  '
Set-TextValue "D21" '14.86'
Set-TextValue "E21" '  -3.17%  '
Set-TextValue "D22" '426.91'
Set-TextValue "E22" '  -1.55%  '
Set-TextValue "D23" '0.594'
Set-TextValue "E23" '  -2.50%  '
Set-TextValue "D24" '77.97'
Set-TextValue "E24" '  -2.04%  '
Set-TextValue "D25" '3.633.48'
Set-TextValue "E25" '  -0.09%  '
Set-TextValue "E26" '  +0.14%  '
Set-TextValue "D27" '0.0000118'
Set-TextValue "E27" '  -1.31%  '
Set-TextValue "D28" '9.30'
Set-TextValue "E28" '  -5.31%  '
Set-TextValue "D29" '7.94'
Set-TextValue "E29" '  -3.90%  '
Set-TextValue "D30" '2.47'
Set-TextValue "E30" '  -0.90%  '
Set-TextValue "E31" '  +0.98%  '
Set-TextValue "D32" '0.166'
Set-TextValue "E32" '  -0.19%  '
Set-TextValue "E33" '  -8.47%  '
Set-TextValue "D34" '25.13'
Set-TextValue "E34" '  -1.00%  '
Set-TextValue "D35" '3.480.22'
Set-TextValue "E35" '  -0.33%  '
Set-TextValue "E37" '  -3.13%  '
Set-TextValue "D38" '5.64'
Set-TextValue "E38" '  -4.66%  '
Set-TextValue "D39" '7.74'
Set-TextValue "E39" '  -3.30%  '
Set-TextValue "E40" '  -0.02%  '
Set-TextValue "D41" '170.07'
Set-TextValue "E41" '  -0.08%  '
Set-TextValue "D42" '0.0860'
Set-TextValue "E42" '  -3.58%  '
Set-TextValue "D43" '5.17'
Set-TextValue "E43" '  -4.77%  '
Set-TextValue "E44" '  -1.69%  '
Set-TextValue "E45" '  -8.66%  '
Set-TextValue "D46" '45.43'
Set-TextValue "E46" '  -0.88%  '
Set-TextValue "D47" '25.95'
Set-TextValue "E47" '  -8.77%  '
Set-TextValue "D48" '1.21'
Set-TextValue "E48" '  -8.16%  '
Set-TextValue "D49" '2.42'
Set-TextValue "E49" '  -0.37%  '
Set-TextValue "D50" '7.15'
Set-TextValue "E50" '  -4.22%  '
Set-TextValue "D51" '0.944'
Set-TextValue "E51" '  -2.82%  '
